$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-11-15 01:18:02"

# Remove all existing hyperlinks first so we can rebuild them cleanly in the correct order/targets
$ws.Hyperlinks.Delete()

$rows = @(
    @{Row=2; B="建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"; D="200,000 円 ~ 300,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434128"; G=368; H="🔥AI,Ai ◆開発"};
    @{Row=3; B="企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)"; D="200,000 円 ~ 300,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434363"; G=348; H="🔥AI,Ai ◆コンサル"};
    @{Row=4; B="【GAS】Yahoo!ショッピング注文完了メール (Gmail) からスプレッドシートに転記する仕事"; D="20,000 円 ~ 50,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5433649"; G=298; H="🔥AI,Ai"};
    @{Row=5; B="画像処理システム(ツール)の開発"; D="20,000 円 ~ 50,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434134"; G=138; H="◆ツール,開発"};
    @{Row=6; B="【システム開発】FileMaker Proを活用した販売システム構築"; D="100,000 円 ~ 200,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434428"; G=118; H="◆開発,システム開発"};
    @{Row=7; B="英語教育の公式LINEアカウント開発・運用スタッフ募集【即日〜3月/4ヶ月/継続可能】"; D="20,000 円 ~ 50,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5433668"; G=63; H="◆開発"};
    @{Row=8; B="【急募】WordPressでの商品検索サイト構築依頼"; D="200,000 円 ~ 300,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5433985"; G=58; H="◇サイト ○WordPress"};
    @{Row=9; B="【GAS活用】業務改善システムの構築依頼"; D="50,000 円 ~ 100,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434156"; G=53; H="◇業務改善"};
    @{Row=10; B="Flutter iOSアプリにおけるRevenueCat導入のバグ修正依頼"; D="50,000 円 ~ 100,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434437"; G=38; H="◇アプリ"};
    @{Row=11; B="wordpressレンダリングを妨げるリソースの除外"; D="200,000 円 ~ 300,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5016989"; G=33; H="○WordPress"};
    @{Row=12; B="【相談から】Laravel7からLaravel12へのサーバーアップデート依頼"; D="500,000 円 ~ 1,000,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5433727"; G=25; H=$null};
    @{Row=13; B="URL付きPDF資料の閲覧状況を可視化し、トラッキングする"; D="100,000 円 ~ 200,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434431"; G=18; H=$null};
    @{Row=14; B="【急募】リマーケティング運用の計測・オーディエンス整備依頼"; D="50,000 円 ~ 100,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434413"; G=18; H=$null};
    @{Row=15; B="初回 ★社内の音響設計スキル向上のため、Modeler / EASE Focus を教えていただける方"; D="50,000 円 ~ 100,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5433823"; G=18; H=$null};
    @{Row=16; B="GAS構築できる方求む!"; D="20,000 円 ~ 50,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434226"; G=13; H=$null};
    @{Row=17; B="月1~5万円以内の小規模タスク依頼"; D="20,000 円 ~ 50,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5433937"; G=13; H=$null};
    @{Row=18; B="【相談のみ】Unityで自動ルート設計プログラムが実現可能か専門家に相談がしたい"; D="1,000 ~ 5,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5434061"; G=10; H=$null};
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $timestamp
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = "システム開発"
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = "期限情報なし"
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    if ($r.H) {
        $ws.Cells.Item($r.Row, 8).Value = $r.H
    } else {
        $ws.Cells.Item($r.Row, 8).Value = ""
    }
    $ws.Hyperlinks.Add($ws.Cells.Item($r.Row, 6), $r.F)
    $ws.Cells.Item($r.Row, 6).Style = "Hyperlink"
}

Write-Host "Done updating rows"
